# Automatische test-sync: 2025-07-23 22:20:50
# Adds a new test-mail log entry (row 13) to the "Logs" sheet, extends the
# conditional-formatting ranges to cover the new row, and updates the
# corresponding tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")

$row = 13

$logs.Cells.Item($row, 1).Value  = "Ik heb maat M ontvangen maar ik had L besteld. Kan ik ruilen?"
$logs.Cells.Item($row, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value  = "Testmail #3: Ik heb maat M ontvangen maar ik had L besteld. Kan ik ruilen?"
$logs.Cells.Item($row, 4).Value  = "Retour / Terugbetaling"
$logs.Cells.Item($row, 5).Value  = "Beste klant,`nBedankt voor je e-mail. Het spijt me te horen dat je de verkeerde maat hebt ontvangen. Om je zo goed mogelijk van dienst te zijn, wil ik vragen of je ons je bestelnummer en de juiste maat kunt doorgeven? Op die manier kunnen we de juiste maat voor je regelen en het omruilproces in gang zetten.`nWe kijken uit naar je reactie en helpen je graag verder.`nMet vriendelijke groet,`n[Naam van het bedrijf]"
$logs.Cells.Item($row, 6).Value  = "2025-07-23 22:20:00"
$logs.Cells.Item($row, 7).Value  = "Ja"
$logs.Cells.Item($row, 8).Value  = "Nee"
$logs.Cells.Item($row, 9).Value  = "Ja"
$logs.Cells.Item($row, 10).Value = "Nee"

# Extend the existing conditional-formatting rules (D, G, H, I, J columns)
# so their range covers the newly added row 13 as well, keeping every rule
# (priority / dxfId / formula) untouched.
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $logs.Range($col + "2:" + $col + "12")
    $newRange = $logs.Range($col + "2:" + $col + "13")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard tally: the "Retour / Terugbetaling" category now has
# two entries instead of one.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 2
